$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 15
$ws.Range("B4").Value = "teacher10"
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = "6->Arabic"
$ws.Range("E4").Value = 1010101010
$ws.Range("F4").Value = "teacher10@"
$ws.Range("G4").Value = " cairo street"
$ws.Range("H4").Value = $false

$ws.Range("A5").Value = 20
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = "1->histroy"
$ws.Range("E5").Value = 2020202020
$ws.Range("F5").Value = "teacher5@gmail.com"
$ws.Range("G5").Value = "street "
$ws.Range("B5").Value = "teacher5"
$ws.Range("H5").Value = $false

[void]$ws.Rows.Item(5).Select()
